$wb = $excel.ActiveWorkbook

# ALC row 64
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H64").Value = 2730.2942
$ws.Range("I64").Value = 3001.5
$ws.Range("J64").Value = 2646.8462
$ws.Range("K64").Value = 3001.5
$ws.Range("L64").Value = 2646.8462
$ws.Range("M64").Value = -2753.5
$ws.Range("N64").Value = -3142.8462

# ALC row 67
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H67").Value = 2730.2942
$ws.Range("I67").Value = 3001.5
$ws.Range("J67").Value = 2646.8462
$ws.Range("K67").Value = 3001.5
$ws.Range("L67").Value = 2646.8462
$ws.Range("M67").Value = -2143.5
$ws.Range("N67").Value = -4362.8462

# ALC row 116
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H116").Value = 1986338.6
$ws.Range("I116").Value = 2316933.5
$ws.Range("J116").Value = 2768.6667
$ws.Range("K116").Value = 2316933.5
$ws.Range("L116").Value = 2768.6667
$ws.Range("M116").Value = -2313491.5
$ws.Range("N116").Value = -9652.6667

# ALC row 128
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H128").Value = 9160.909
$ws.Range("I128").Value = 0
$ws.Range("J128").Value = 9160.909
$ws.Range("K128").Value = 0
$ws.Range("L128").Value = 9160.909
$ws.Range("N128").Value = -19120.909

# ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2688.4517
$ws.Range("I132").Value = 2611
$ws.Range("J132").Value = 3211.25
$ws.Range("K132").Value = 7833
$ws.Range("L132").Value = 9633.75
$ws.Range("M132").Value = -5303
$ws.Range("N132").Value = -14693.75

# ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2825.2415
$ws.Range("I138").Value = 1474.4667
$ws.Range("J138").Value = 3296.442
$ws.Range("K138").Value = 4423.4001
$ws.Range("L138").Value = 9889.326000000001
$ws.Range("M138").Value = 716.5999000000002
$ws.Range("N138").Value = -20169.326

# ARM row 2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 773.38464
$ws.Range("I2").Value = 795.3333
$ws.Range("J2").Value = 754.5714
$ws.Range("K2").Value = 795.3333
$ws.Range("L2").Value = 754.5714
$ws.Range("M2").Value = -682.3333
$ws.Range("N2").Value = -980.5714

# ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 12235.87
$ws.Range("I32").Value = 3840.6616
$ws.Range("J32").Value = 27826.97
$ws.Range("K32").Value = 3840.6616
$ws.Range("L32").Value = 27826.97
$ws.Range("M32").Value = -3553.6616
$ws.Range("N32").Value = -28400.97

# ARM row 45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 2402.6924
$ws.Range("I45").Value = 2754.375
$ws.Range("J45").Value = 1840
$ws.Range("K45").Value = 2754.375
$ws.Range("L45").Value = 1840
$ws.Range("M45").Value = -2377.375
$ws.Range("N45").Value = -2594

# ARM row 102
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 2500
$ws.Range("I102").Value = 2500
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 2500
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -878

# ARM row 116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 773.38464
$ws.Range("I116").Value = 795.3333
$ws.Range("J116").Value = 754.5714
$ws.Range("K116").Value = 795.3333
$ws.Range("L116").Value = 754.5714
$ws.Range("M116").Value = 1498.6667
$ws.Range("N116").Value = -5342.5714

# ARM row 121
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H121").Value = 0
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()

# ARM row 122
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H122").Value = 1919.1724
$ws.Range("I122").Value = 1468.8572
$ws.Range("J122").Value = 3101.25
$ws.Range("K122").Value = 4406.571599999999
$ws.Range("L122").Value = 9303.75
$ws.Range("M122").Value = -1956.571599999999
$ws.Range("N122").Value = -14203.75

# BSM row 3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 773.38464
$ws.Range("I3").Value = 795.3333
$ws.Range("J3").Value = 754.5714
$ws.Range("K3").Value = 795.3333
$ws.Range("L3").Value = 754.5714
$ws.Range("M3").Value = -681.3333
$ws.Range("N3").Value = -982.5714

# BSM row 99
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1307.48
$ws.Range("I99").Value = 1043.7778
$ws.Range("J99").Value = 1985.5714
$ws.Range("K99").Value = 1043.7778
$ws.Range("L99").Value = 1985.5714
$ws.Range("M99").Value = 454.2221999999999
$ws.Range("N99").Value = -4981.5714

# BSM row 103
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H103").Value = 30000
$ws.Range("I103").Value = 0
$ws.Range("J103").Value = 30000
$ws.Range("K103").Value = 0
$ws.Range("L103").Value = 30000
$ws.Range("N103").Value = -32344

# BSM row 105
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1801.6666
$ws.Range("I105").Value = 1003.3333
$ws.Range("J105").Value = 2600
$ws.Range("K105").Value = 1003.3333
$ws.Range("L105").Value = 2600
$ws.Range("M105").Value = 743.6667

# BSM row 107
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 1346.16
$ws.Range("I107").Value = 1124.7646
$ws.Range("J107").Value = 1816.625
$ws.Range("K107").Value = 1124.7646
$ws.Range("L107").Value = 1816.625
$ws.Range("M107").Value = 795.2354
$ws.Range("N107").Value = -5656.625

# CRP row 58
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 1364.8462
$ws.Range("I58").Value = 976.2381
$ws.Range("J58").Value = 1818.2222
$ws.Range("K58").Value = 976.2381
$ws.Range("L58").Value = 1818.2222
$ws.Range("M58").Value = -773.2381
$ws.Range("N58").Value = -2224.2222

# CRP row 94
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H94").Value = 1000
$ws.Range("I94").Value = 0
$ws.Range("J94").Value = 1000
$ws.Range("K94").Value = 0
$ws.Range("L94").Value = 1000
$ws.Range("M94").ClearContents()
$ws.Range("N94").Value = -1902

# CRP row 134
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 4902.9165
$ws.Range("I134").Value = 5466.7334
$ws.Range("J134").Value = 3963.2222
$ws.Range("K134").Value = 16400.2002
$ws.Range("L134").Value = 11889.6666
$ws.Range("M134").Value = -13865.2002
$ws.Range("N134").Value = -16959.6666

# CRP row 136
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H136").Value = 1364.8462
$ws.Range("I136").Value = 976.2381
$ws.Range("J136").Value = 1818.2222
$ws.Range("K136").Value = 2928.7143
$ws.Range("L136").Value = 5454.6666
$ws.Range("M136").Value = -378.7143000000001
$ws.Range("N136").Value = -10554.6666

# CUL row 131
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1019.53656
$ws.Range("I131").Value = 622.75
$ws.Range("J131").Value = 1062.4324
$ws.Range("K131").Value = 1868.25
$ws.Range("L131").Value = 3187.2972
$ws.Range("M131").Value = 3171.75
$ws.Range("N131").Value = -13267.2972

# CUL row 134
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H134").Value = 2393.889
$ws.Range("I134").Value = 3269
$ws.Range("J134").Value = 1300
$ws.Range("K134").Value = 9807
$ws.Range("L134").Value = 3900
$ws.Range("M134").Value = -4737

# GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 911.26666
$ws.Range("I102").Value = 896.7778
$ws.Range("J102").Value = 933
$ws.Range("K102").Value = 896.7778
$ws.Range("L102").Value = 933
$ws.Range("M102").Value = 725.2222
$ws.Range("N102").Value = -4177

# GSM row 107
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 774
$ws.Range("I107").Value = 285.875
$ws.Range("J107").Value = 1262.125
$ws.Range("K107").Value = 285.875
$ws.Range("L107").Value = 1262.125
$ws.Range("M107").Value = 1634.125
$ws.Range("N107").Value = -5102.125

# GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 2831.4
$ws.Range("I122").Value = 4052.3333
$ws.Range("J122").Value = 1000
$ws.Range("K122").Value = 12156.9999
$ws.Range("L122").Value = 3000
$ws.Range("M122").Value = -9706.999899999999
$ws.Range("N122").Value = -7900

# LTW row 82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 1951.125
$ws.Range("I82").Value = 1184.8334
$ws.Range("J82").Value = 4250
$ws.Range("K82").Value = 1184.8334
$ws.Range("L82").Value = 4250
$ws.Range("M82").Value = -823.8334
$ws.Range("N82").Value = -4972

# LTW row 85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 1951.125
$ws.Range("I85").Value = 1184.8334
$ws.Range("J85").Value = 4250
$ws.Range("K85").Value = 1184.8334
$ws.Range("L85").Value = 4250
$ws.Range("M85").Value = 63.16660000000002
$ws.Range("N85").Value = -6746

# LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 1531.6666
$ws.Range("I100").Value = 1218
$ws.Range("J100").Value = 1923.75
$ws.Range("K100").Value = 1218
$ws.Range("L100").Value = 1923.75
$ws.Range("M100").Value = -677

# WVR row 81
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 45456724
$ws.Range("I81").Value = 90911224
$ws.Range("J81").Value = 2218.5454
$ws.Range("K81").Value = 181822448
$ws.Range("L81").Value = 4437.0908
$ws.Range("M81").Value = -181821387
$ws.Range("N81").Value = -6559.0908

# WVR row 84
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H84").Value = 45456724
$ws.Range("I84").Value = 90911224
$ws.Range("J84").Value = 2218.5454
$ws.Range("K84").Value = 909112240
$ws.Range("L84").Value = 22185.454
$ws.Range("M84").Value = -909106936
$ws.Range("N84").Value = -32793.454
